# Re-pull data / push all data / mean calculation
# Update the "dSF" (F) column with refreshed computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    4  = 1
    5  = 1
    6  = 4
    9  = 2
    10 = 4
    12 = -1
    13 = -1
    14 = 2
    15 = 1
    16 = -1
    17 = 2
    18 = 1
    19 = -4
    20 = 2
    21 = 3
    22 = -5
    23 = -5
    24 = -3
    25 = 1
    26 = -1
    27 = 2
    28 = -4
    29 = 3
    30 = 2
    31 = -1
    33 = -2
    34 = -4
    35 = -3
    36 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
